$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: was "Interactive Binning" / "Show initial (after) summary statistics table panel"
# becomes "Preview & Download Settings" / "Show preview dataset", with new Status/Date
$ws.Range("B29").Value = "Preview & Download Settings"
$ws.Range("C29").Value = "Show preview dataset"
$ws.Range("D29").Value = "-"
$ws.Range("E29").Value = "Closed"
$ws.Range("F29").Value = (Get-Date -Year 2023 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Row 30: was "Interactive Binning" (only B filled)
# becomes "Preview & Download Settings" / "Download json bin settings button"
$ws.Range("B30").Value = "Preview & Download Settings"
$ws.Range("C30").Value = "Download json bin settings button"
$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = "Closed"
$ws.Range("F30").Value = (Get-Date -Year 2023 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Row 31: was "Interactive Binning" (only B filled)
# becomes "(Backend)" / "Initial binning & save result to shared storage - treat every value as a unique bin"
$ws.Range("B31").Value = "(Backend)"
$ws.Range("C31").Value = "Initial binning & save result to shared storage - treat every value as a unique bin"
$ws.Range("D31").Value = "-"
$ws.Range("E31").Value = "Closed"
$ws.Range("F31").Value = (Get-Date -Year 2023 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Widen column B to fit new, longer strings (closest achievable snap to target 25.453125)
$ws.Columns.Item(2).ColumnWidth = 24.6

# Update view: scroll down a bit and move selection to F32
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F32").Select()
